# update pin allocations for I2C OLED
#
# Applies the six table-cell edits described by the commit:
#   1. TxD row (pin 4 / GPIO.7 .. pin 14 / TxD): last (Role) cell gains the
#      "_GoBack" bookmark that Word drops at the most-recent-edit location.
#   2. MOSI row (pin 10): the "SPI OLED" label is cleared.
#   3. MISO row (pin 9 .. pin 25 / GPIO. 6): the "_GoBack" bookmark that used
#      to sit on the "GPIO. 6" cell is removed (text is kept).
#   4. SCLK row (pin 11): the "SPI Clock OLED" label is cleared.
#   5. SDA.0 / SCL.0 row: the two previously-empty Role cells are filled in
#      with "I2C Data OLED" and "I2C Clock OLED" respectively.
#   6. Last row (pin 21 / GPIO.29): "OLED " + "DC Pin" runs are normalized
#      into a single "OLED DC Pin" run.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. TxD row: add the "_GoBack" bookmark to the (empty) last cell -------
$cell = $t.Cell(5, 7)
$d.Bookmarks.Add("_GoBack", $cell.Range)

# --- 2. MOSI row: remove the "SPI OLED" text, leaving the cell empty -------
$t = $d.Tables.Item(1)
$cell = $t.Cell(11, 3)
$cell.Range.Find.Execute("SPI OLED", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 2)

# --- 3. MISO row: drop the stray "_GoBack" bookmark, keep "GPIO. 6" --------
$t = $d.Tables.Item(1)
$cell = $t.Cell(12, 6)
$cell.Range.Delete()
$t = $d.Tables.Item(1)
$cell = $t.Cell(12, 6)
$cell.Range.InsertAfter("GPIO. 6")

# --- 4. SCLK row: remove the "SPI Clock OLED" text -------------------------
$t = $d.Tables.Item(1)
$cell = $t.Cell(13, 3)
$cell.Range.Find.Execute("SPI Clock OLED", $false, $false, $false, $false, `
    $false, $true, 1, $false, "", 2)

# --- 5. SDA.0 / SCL.0 row: fill in the two new labels -----------------------
$t = $d.Tables.Item(1)
$cell = $t.Cell(15, 3)
$cell.Range.Text = "I2C Data OLED"

$t = $d.Tables.Item(1)
$cell = $t.Cell(15, 7)
$cell.Range.Text = "I2C Clock OLED"

# --- 6. Last row: merge "OLED " + "DC Pin" into a single run ---------------
$t = $d.Tables.Item(1)
$cell = $t.Cell(21, 7)
$cell.Range.Find.Execute("OLED DC Pin", $false, $false, $false, $false, `
    $false, $true, 1, $false, "OLED DC Pin", 2)
